$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.2612
$ws.Range("C7").Value = -12.80810000000001
$ws.Range("B9").Value = 5.113099999999998
$ws.Range("C12").Value = -10.865
$ws.Range("B13").Value = 6.377199999999998
$ws.Range("C14").Value = -11.9312
$ws.Range("E15").Value = 16.26280000000001
$ws.Range("B16").Value = 4.635000000000002
$ws.Range("B18").Value = 6.801599999999997
$ws.Range("C19").Value = -11.63350000000001
$ws.Range("B20").Value = 8.877599999999997
$ws.Range("B26").Value = 4.394700000000003
$ws.Range("C26").Value = -12.6874
$ws.Range("B27").Value = 5.401000000000002
$ws.Range("C27").Value = -12.54919999999999
$ws.Range("E28").Value = 16.18989999999999
$ws.Range("B29").Value = 4.728499999999998
$ws.Range("C29").Value = -10.58810000000001
$ws.Range("E33").Value = 17.05470000000001
$ws.Range("B35").Value = 8.649699999999998
$ws.Range("E35").Value = 16.71
$ws.Range("B36").Value = 9.172200000000002
$ws.Range("C37").Value = -12.4782
$ws.Range("C38").Value = -12.6626
$ws.Range("E38").Value = 16.43119999999999
$ws.Range("E43").Value = 17.17990000000001
$ws.Range("E44").Value = 16.65659999999999
$ws.Range("B45").Value = 5.669800000000002
$ws.Range("E45").Value = 16.58440000000001
$ws.Range("C47").Value = -12.6058
$ws.Range("E47").Value = 16.15019999999999
$ws.Range("C51").Value = -11.2442
$ws.Range("E51").Value = 17.29330000000001
$ws.Range("C52").Value = -11.3805
$ws.Range("E54").Value = 16.5122
$ws.Range("B55").Value = 6.256299999999998
$ws.Range("C55").Value = -13.90160000000001
$ws.Range("B57").Value = 5.262399999999993
$ws.Range("E57").Value = 16.2275
$ws.Range("E62").Value = 16.2471
$ws.Range("E63").Value = 18.75610000000001
$ws.Range("E67").Value = 17.20520000000002
$ws.Range("B69").Value = 5.830899999999994
$ws.Range("C69").Value = -11.3887
$ws.Range("C70").Value = -11.3257
$ws.Range("E70").Value = 17.43480000000002
$ws.Range("B76").Value = 5.652100000000005
$ws.Range("C76").Value = -12.023
$ws.Range("B78").Value = 9.823900000000004
$ws.Range("C81").Value = -13.26229999999999
$ws.Range("E81").Value = 16.8013
$ws.Range("B82").Value = 5.528900000000002
$ws.Range("B83").Value = 5.658999999999996
$ws.Range("C83").Value = -14.0247
$ws.Range("E88").Value = 16.2444
$ws.Range("B93").Value = 5.619999999999997
$ws.Range("C94").Value = -10.24670000000001
$ws.Range("E96").Value = 16.35549999999999
$ws.Range("B97").Value = 5.524899999999998
$ws.Range("E99").Value = 16.58010000000001
$ws.Range("C100").Value = -12.37969999999999
$ws.Range("C102").Value = -12.9943
